# Update column G ("K" = strikeouts) values for rows 2-30 in Sheet1
# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @{
    2  = 0
    3  = 3
    4  = 6
    5  = 3
    6  = 3
    7  = 2
    8  = 3
    9  = 4
    10 = 2
    11 = 4
    12 = 5
    13 = 3
    14 = 6
    15 = 1
    16 = 3
    17 = 7
    18 = 2
    19 = 1
    20 = 2
    21 = 9
    22 = 4
    23 = 5
    24 = 5
    25 = 5
    26 = 4
    27 = 3
    28 = 5
    29 = 2
    30 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}

$wb.Save()
